$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 193, shifting existing rows 193:204 down to 194:205.
$ws.Rows(193).Insert()

# Populate the newly inserted row 193 with the new weekly price observation.
$ws.Cells.Item(193, 1).Value = 9
$ws.Cells.Item(193, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(193, 3).Value = "Metropolitana"
$ws.Cells.Item(193, 4).Value = "1/17/2022"
$ws.Cells.Item(193, 5).Value = 13
$ws.Cells.Item(193, 6).Value = 100112030
$ws.Cells.Item(193, 7).Value = "Poroto granado"
$ws.Cells.Item(193, 8).Value = "Sin especificar"
$ws.Cells.Item(193, 9).Value = "Primera"
$ws.Cells.Item(193, 10).Value = 61
$ws.Cells.Item(193, 11).Value = 25000
$ws.Cells.Item(193, 12).Value = 28000
$ws.Cells.Item(193, 13).Value = 26525
$ws.Cells.Item(193, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(193, 15).Value = "Región del Maule"
$ws.Cells.Item(193, 16).Value = 1061
$ws.Cells.Item(193, 17).Value = 25
$ws.Cells.Item(193, 18).Value = "Hortaliza"
